$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "password" column (E) values: strip the md5(...) wrapper,
# and rename the last student's placeholder from ubahyah0 -> ubahyah10 ---
$ws.Range("E2").Value  = "ubahyah1"
$ws.Range("E3").Value  = "ubahyah2"
$ws.Range("E4").Value  = "ubahyah3"
$ws.Range("E5").Value  = "ubahyah4"
$ws.Range("E6").Value  = "ubahyah5"
$ws.Range("E7").Value  = "ubahyah6"
$ws.Range("E8").Value  = "ubahyah7"
$ws.Range("E9").Value  = "ubahyah8"
$ws.Range("E10").Value = "ubahyah9"
$ws.Range("E11").Value = "ubahyah10"

# --- Remove the last six member rows entirely (now blank) ---
$ws.Range("A12:E17").ClearContents()

# --- Freeze the header row and leave the selection on E4 ---
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
[void]$ws.Range("E4").Select()

# --- Page setup: portrait orientation ---
[void]($ws.PageSetup.Orientation = 1)
